$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(5730,5733,5665,5673,5580,5589,5556,5612,5527,5512,5488,5506,5579,5536,5566,5590,5692,5698,5790,5879,6008,6131,6228,6403,6685,6825,6924,6980,7055,7078,7098,7009,6976,6911,6816,6639,6585,6476,6404,6288,6175,6149,6072,6024,5941,5852,5834,5851,5892,5816,5770,5823,5853,5826,5779,5824,5805,5799,5838,5912,5950,5949,6020,6077,6269)
$bValues = @(45751,45751.01041666666,45751.02083333334,45751.03125,45751.04166666666,45751.05208333334,45751.0625,45751.07291666666,45751.08333333334,45751.09375,45751.10416666666,45751.11458333334,45751.125,45751.13541666666,45751.14583333334,45751.15625,45751.16666666666,45751.17708333334,45751.1875,45751.19791666666,45751.20833333334,45751.21875,45751.22916666666,45751.23958333334,45751.25,45751.26041666666,45751.27083333334,45751.28125,45751.29166666666,45751.30208333334,45751.3125,45751.32291666666,45751.33333333334,45751.34375,45751.35416666666,45751.36458333334,45751.375,45751.38541666666,45751.39583333334,45751.40625,45751.41666666666,45751.42708333334,45751.4375,45751.44791666666,45751.45833333334,45751.46875,45751.47916666666,45751.48958333334,45751.5,45751.51041666666,45751.52083333334,45751.53125,45751.54166666666,45751.55208333334,45751.5625,45751.57291666666,45751.58333333334,45751.59375,45751.60416666666,45751.61458333334,45751.625,45751.63541666666,45751.64583333334,45751.65625,45751.66666666666)

$n = $aValues.Length

for ($i = 0; $i -lt $n; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Ensure the newly added rows in column B carry the same date/time
# number format as the pre-existing rows (style index 2 -> "YYYY-MM-DD HH:MM:SS").
$ws.Range("B2:B" + (1 + $n)).NumberFormat = "YYYY-MM-DD HH:MM:SS"
